$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-15 16:43:33"

$wsZhCn.Range("H4").Value = "2016-08-15 16:43:28"
$wsZhCn.Range("K4").Value = "2016-08-15 16:43:45"

$wsDeDe.Range("H4").Value = "2016-08-15 16:43:33"
$wsDeDe.Range("K4").Value = "2016-08-15 16:43:52"
